$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "vendor"
$ws.Range("C3").Value = "Cisco"
$ws.Range("C2").Value = "Juniper"
$ws.Range("C4").Value = "Nokia"

$ws.Range("C5").Select()
